# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" everywhere it
#   is used (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - The per-locale "Latest Handoff Datetime" / "Latest HO Xliff Generate
#   Date" timestamps are bumped to the new handoff-generation time.
# - The Status/locale columns widen slightly to fit the new, longer status
#   text ("Ready for handoff" vs "In Translation").

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps bumped to reflect the new handoff generation ---
$wsZhCn.Range("H2").Value = "2016-08-12 09:11:35"
$wsDeDe.Range("H2").Value = "2016-08-12 09:11:42"
$wsOverview.Range("G2").Value = "2016-08-12 09:11:42"

# --- Widen the status/locale columns to fit the new text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
